$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "isophonics_46"
$ws.Range("B2").Value = "schubert-winterreise_129"
$ws.Range("C2").Value = 0.1405405405405405
$ws.Range("D2").Value = "[['D', 'G', 'D']]"
$ws.Range("E2").Value = "[['F:maj', 'A#:maj', 'F:maj']]"
$ws.Range("F2").Value = "[(28.921995, 33.101587)]"
$ws.Range("G2").Value = "[(105.0, 108.24)]"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "spotify:track:1nvxQGWCnikMK7a4HYQvSx"

# Row 3
$ws.Range("A3").Value = "isophonics_171"
$ws.Range("B3").Value = "isophonics_150"
$ws.Range("C3").Value = 0.1471291866028708
$ws.Range("D3").Value = "[['E', 'E', 'C#:min']]"
$ws.Range("E3").Value = "[['G', 'G', 'E:min']]"
$ws.Range("F3").Value = "[(24.616991, 31.687468)]"
$ws.Range("G3").Value = "[(47.639639, 53.885807)]"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# Row 4
$ws.Range("A4").Value = "schubert-winterreise_63"
$ws.Range("B4").Value = "schubert-winterreise_118"
$ws.Range("C4").Value = 0.07984496124031007
$ws.Range("D4").Value = "[['G:7', 'C:min', 'B:dim7/C', 'C:min'], ['G:7', 'C:min', 'C:min', 'D:hdim7/C'], ['C:min', 'B:dim7/C', 'C:min', 'B:dim7/C']]"
$ws.Range("E4").Value = "[['E:7', 'A:min', 'G#:dim7', 'A:min'], ['E:7', 'A:min', 'A:min', 'B:hdim7/D'], ['A:min', 'G#:dim7', 'A:min', 'G#:dim7']]"
$ws.Range("F4").Value = "[(78.74, 89.26), (26.48, 34.32), (0.24, 9.6)]"
$ws.Range("G4").Value = "[(20.38, 24.76), (14.2, 20.38), (21.9, 26.28)]"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "spotify:track:2qCvEz2hEb92VFATqVvrht"

# Row 5
$ws.Range("A5").Value = "isophonics_279"
$ws.Range("B5").Value = "schubert-winterreise_114"
$ws.Range("C5").Value = 0.2259615384615385
$ws.Range("D5").Value = "[['C:maj', 'F/5', 'C:maj']]"
$ws.Range("E5").Value = "[['D:maj/F#', 'G:maj', 'D:maj']]"
$ws.Range("F5").Value = "[(10.619, 16.742)]"
$ws.Range("G5").Value = "[(57.48, 64.58)]"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("A6").Value = "schubert-winterreise_3"
$ws.Range("B6").Value = "jaah_3"
$ws.Range("C6").Value = 0.09027777777777778
$ws.Range("D6").Value = "[['G#:7', 'C#:maj', 'G#:7']]"
$ws.Range("E6").Value = "[['Bb:7', 'Eb', 'Bb:7']]"
$ws.Range("F6").Value = "[(7.34, 9.06)]"
$ws.Range("G6").Value = "[(44.07, 46.31)]"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

# Row 7
$ws.Range("A7").Value = "schubert-winterreise_139"
$ws.Range("B7").Value = "schubert-winterreise_164"
$ws.Range("C7").Value = 0.2153846153846154
$ws.Range("D7").Value = "[['E:maj/B', 'B:7', 'E:maj/B'], ['D#:min/A#', 'A#:7', 'D#:min']]"
$ws.Range("E7").Value = "[['A#:maj/F', 'F:7', 'A#:maj'], ['A:min', 'E:7', 'A:min']]"
$ws.Range("F7").Value = "[(84.96, 91.38), (116.66, 121.42)]"
$ws.Range("G7").Value = "[(114.3, 123.14), (19.82, 29.1)]"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = "spotify:track:3OD2uwEUQKg0WyW9Lewata"

# Row 8
$ws.Range("A8").Value = "schubert-winterreise_77"
$ws.Range("B8").Value = "schubert-winterreise_111"
$ws.Range("C8").Value = 0.0945054945054945
$ws.Range("D8").Value = "[['F:min', 'F#/A#', 'F:min/C'], ['G:hdim7/A#', 'C:7', 'F:min']]"
$ws.Range("E8").Value = "[['G:min', 'G#:maj/C', 'G:min'], ['A:hdim7/D#', 'D:7', 'G:min']]"
$ws.Range("F8").Value = "[(187.92, 190.22), (23.18, 30.14)]"
$ws.Range("G8").Value = "[(42.58, 45.2), (100.38, 113.76)]"
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"

# Row 9
$ws.Range("A9").Value = "schubert-winterreise_141"
$ws.Range("B9").Value = "jaah_9"
$ws.Range("C9").Value = 0.1785714285714285
$ws.Range("D9").Value = "[['C:7', 'F:maj', 'F:maj']]"
$ws.Range("E9").Value = "[['C:7', 'F', 'F']]"
$ws.Range("F9").Value = "[(73.74, 83.38)]"
$ws.Range("G9").Value = "[(51.69, 55.72)]"
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = "schubert-winterreise_170"
$ws.Range("B10").Value = "schubert-winterreise_162"
$ws.Range("C10").Value = 0.2657342657342657
$ws.Range("D10").Value = "[['G:maj', 'D:7/C', 'G:maj/B', 'D:7/C', 'G:maj/B']]"
$ws.Range("E10").Value = "[['G:maj/D', 'D:7', 'G:maj', 'D:7', 'G:maj']]"
$ws.Range("F10").Value = "[(72.52, 90.98)]"
$ws.Range("G10").Value = "[(23.26, 29.8)]"
$ws.Range("H10").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"
$ws.Range("I10").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"

# Row 11
$ws.Range("A11").Value = "schubert-winterreise_105"
$ws.Range("B11").Value = "schubert-winterreise_163"
$ws.Range("C11").Value = 0.2528735632183908
$ws.Range("D11").Value = "[['G:min', 'D:7/G', 'G:min', 'D:maj/G', 'G:min']]"
$ws.Range("E11").Value = "[['F:min', 'C:7', 'F:min', 'C:maj', 'F:min']]"
$ws.Range("F11").Value = "[(23.74, 51.72)]"
$ws.Range("G11").Value = "[(10.58, 17.96)]"
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = "spotify:track:1nvxQGWCnikMK7a4HYQvSx"

# Row 12
$ws.Range("A12").Value = "schubert-winterreise_181"
$ws.Range("B12").Value = "isophonics_157"
$ws.Range("C12").Value = 0.1964285714285714
$ws.Range("D12").Value = "[['D#:maj', 'D#:7', 'G#:maj'], ['D#:maj', 'A#:maj', 'D#:maj']]"
$ws.Range("E12").Value = "[['E', 'E:7/3', 'A'], ['E', 'B', 'E']]"
$ws.Range("F12").Value = "[(17.98, 26.32), (2.56, 21.44)]"
$ws.Range("G12").Value = "[(20.410362, 24.856984), (12.921927, 21.908049)]"
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = "schubert-winterreise_23"
$ws.Range("B13").Value = "schubert-winterreise_63"
$ws.Range("C13").Value = 0.1916666666666667
$ws.Range("D13").Value = "[['G:maj', 'D:7', 'G:maj']]"
$ws.Range("E13").Value = "[['C/G', 'G:7', 'C']]"
$ws.Range("F13").Value = "[(7.56, 16.36)]"
$ws.Range("G13").Value = "[(253.84, 257.6)]"
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("A14").Value = "schubert-winterreise_88"
$ws.Range("B14").Value = "schubert-winterreise_167"
$ws.Range("C14").Value = 0.3314393939393939
$ws.Range("D14").Value = "[['A:maj/E', 'E:7', 'A:maj', 'E:7', 'A:maj']]"
$ws.Range("E14").Value = "[['C:maj', 'G:7', 'C:maj', 'G:7', 'C:maj']]"
$ws.Range("F14").Value = "[(16.04, 21.0)]"
$ws.Range("G14").Value = "[(0.58, 10.18)]"
$ws.Range("H14").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"
$ws.Range("I14").Value = ""

# Row 15
$ws.Range("A15").Value = "isophonics_0"
$ws.Range("B15").Value = "isophonics_290"
$ws.Range("C15").Value = 0.08977272727272727
$ws.Range("D15").Value = "[['Db', 'Gb', 'Ab'], ['Db', 'Ab', 'Db']]"
$ws.Range("E15").Value = "[['C', 'F', 'G'], ['C', 'G', 'C']]"
$ws.Range("F15").Value = "[(0.344657, 5.571955), (62.09653, 72.37136)]"
$ws.Range("G15").Value = "[(25.807392, 30.753242), (8.482392, 11.666439)]"
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = "jaah_21"
$ws.Range("B16").Value = "isophonics_223"
$ws.Range("C16").Value = 0.1340659340659341
$ws.Range("D16").Value = "[['Ab:7', 'Db/3', 'Db:min/b3', 'Ab/b5']]"
$ws.Range("E16").Value = "[['A:7', 'D', 'D:min', 'A']]"
$ws.Range("F16").Value = "[(3.3, 5.85)]"
$ws.Range("G16").Value = "[(12.823786, 18.895804)]"
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = "spotify:track:3KfbEIOC7YIv90FIfNSZpo"

# Row 17
$ws.Range("A17").Value = "isophonics_271"
$ws.Range("B17").Value = "jaah_69"
$ws.Range("C17").Value = 0.1405405405405405
$ws.Range("D17").Value = "[['E', 'E/5', 'E']]"
$ws.Range("E17").Value = "[['Eb', 'Eb', 'Eb']]"
$ws.Range("F17").Value = "[(61.881201, 65.886643)]"
$ws.Range("G17").Value = "[(16.13, 24.98)]"
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
